$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New S-val data (regenerated to filter save games)
$data = @{
    2 = @(0.04271373187048222, 0.04071648406533734, 0.1494219747398047, 0.4942365360607697, 1, 0.7270887267363939)
    3 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 4.358119930609447)
    4 = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 3.56341032713086)
    5 = @(0.01293466051926884, 0.002571899574220771, 3.537761648806719, 0.4942365360607697, 0, 4.047504744960978)
    6 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    7 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
    8 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 8.974608811992548)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("F$r").Value = $vals[4]
    $ws.Range("G$r").Value = $vals[5]
}
